# Fruta / hortaliza, semanal
# Insert a new weekly record at row 202 (Femacal de La Calera - Chirimoya,
# Provincia del Elquí, 2022-10-21, Primera) and push the remaining existing
# rows (old 202-252) down by one, extending the table to row 253.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(202).Insert()

$ws.Range("A202").Value() = 3
$ws.Range("B202").Value() = "Femacal de La Calera"
$ws.Range("C202").Value() = "Coquimbo"
$ws.Range("D202").Value() = 44855
$ws.Range("E202").Value() = 5
$ws.Range("F202").Value() = "Fruta"
$ws.Range("G202").Value() = 100107
$ws.Range("H202").Value() = "Otros"
$ws.Range("I202").Value() = 100107002
$ws.Range("J202").Value() = "Chirimoya"
$ws.Range("K202").Value() = "Cultivar IV Región"
$ws.Range("L202").Value() = "Primera"
$ws.Range("M202").Value() = 45
$ws.Range("N202").Value() = 27000
$ws.Range("O202").Value() = 27000
$ws.Range("P202").Value() = 27000
$ws.Range("Q202").Value() = "$/bandeja 10 kilos"
$ws.Range("R202").Value() = "Provincia del Elquí"
$ws.Range("S202").Value() = 2700
$ws.Range("T202").Value() = 10
